$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A for rows 2-13: monthly_gross_earnings_effect_* -> yearly_gross_earnings_effect_*
$ws.Range("A2").Value = "yearly_gross_earnings_effect_female_west_year1"
$ws.Range("A3").Value = "yearly_gross_earnings_effect_female_east_year1"
$ws.Range("A4").Value = "yearly_gross_earnings_effect_male_west_year1"
$ws.Range("A5").Value = "yearly_gross_earnings_effect_male_east_year1"
$ws.Range("A6").Value = "yearly_gross_earnings_effect_female_west_year2"
$ws.Range("A7").Value = "yearly_gross_earnings_effect_female_east_year2"
$ws.Range("A8").Value = "yearly_gross_earnings_effect_male_west_year2"
$ws.Range("A9").Value = "yearly_gross_earnings_effect_male_east_year2"
$ws.Range("A10").Value = "yearly_gross_earnings_effect_female_west_year3"
$ws.Range("A11").Value = "yearly_gross_earnings_effect_female_east_year3"
$ws.Range("A12").Value = "yearly_gross_earnings_effect_male_west_year3"
$ws.Range("A13").Value = "yearly_gross_earnings_effect_male_east_year3"

# Update the selection to match the post-edit state (cursor moved to A27)
$ws.Range("A27").Select()
